$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values following Dr Hou's advice (Ligand/Receptor-expressing cell counts
# changed from 1 to 3, with downstream statistics recomputed accordingly).

$data = @{
    2 = @{ E=3; G=4.278141666666667;  H=12.834425;  I=0.9663225094340192; J=0.9663225094340191; K=3; M=28.86405866666667; N=86.59217600000001; O=0.3025739760541936; P=0.3025739760541936; Q=123.4845320509778; R=1111.3607884588;  S=0.2923840438301172; T=0.2923840438301171 }
    3 = @{ E=3; G=4.278141666666667;  H=12.834425;  I=0.9663225094340192; J=0.9663225094340191; K=3; M=42.68037399999999; N=128.041122;        O=0.4474066038250392; P=0.4474066038250392; Q=182.5926863583166; R=1643.33417722485; S=0.4323390721455639; T=0.4323390721455638 }
    4 = @{ E=3; G=4.278141666666667;  H=12.834425;  I=0.9663225094340192; J=0.9663225094340191; K=3; M=23.85061433333334; N=71.55184300000001; O=0.2500194201207672; P=0.2500194201207672; Q=102.0363069550306; R=918.326762595275;  S=0.2415993934583381; T=0.2415993934583381 }
    5 = @{ E=3; G=0.1490983333333333; H=0.447295;   I=0.03367749056598092; J=0.03367749056598091; K=3; M=28.86405866666667; N=86.59217600000001; O=0.3025739760541936; P=0.3025739760541936; Q=4.303583040435556; R=38.73224736392;    S=0.01018993222407644; T=0.01018993222407644 }
    6 = @{ E=3; G=0.1490983333333333; H=0.447295;   I=0.03367749056598092; J=0.03367749056598091; K=3; M=42.68037399999999; N=128.041122;        O=0.4474066038250392; P=0.4474066038250392; Q=6.363572629443333; R=57.27215366498999; S=0.01506753167947532; T=0.01506753167947532 }
    7 = @{ E=3; G=0.1490983333333333; H=0.447295;   I=0.03367749056598092; J=0.03367749056598091; K=3; M=23.85061433333334; N=71.55184300000001; O=0.2500194201207672; P=0.2500194201207672; Q=3.556086846076111; R=32.004781614685;  S=0.008420026662429159; T=0.008420026662429157 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
